$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1489
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202402/hdlmhoLp1708932790894.jpeg"
    $ws.Range("F5").Value = 2271
    $ws.Range("F7").Value = 1392
    $ws.Range("F9").Value = 149
    $ws.Range("F11").Value = 331
    $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202402/9cMJMElF1708938074308.png"
}
